# ARKCORR-18 Use regular LDAP groups as owning group, versus ad-hoc groups.
# Uppercase the LDAP group names used in the "Set Owning Group" action column (E)
# of the "Save Case File Rules" rule table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E31").Value = "owning group, new String('ANALYST - AG')"
$ws.Range("E33").Value = "owning group, new String('EXECSEC APPROVAL')"
$ws.Range("E34").Value = "owning group, new String('EXECSEC RELEASE')"
$ws.Range("E30").Value = "owning group, new String('EXECSEC INTAKE')"
$ws.Range("E32").Value = "owning group, new String('SUPERVISOR - AG')"

# Reflect the editor's final on-screen view state (scrolled down, cell E33 selected)
$ws.Activate()
$ws.Range("E33").Select()
